# Applies the September courtesy-notice edit:
#   1. Bump the letter date from September 19 -> September 21, 2025.
#   2. Split the single "929 Story Road, San Jose CA 95122" mailing-address
#      line into two separate paragraphs ("929 Story Road" / "San Jose, CA
#      95122"), matching the new mailing-address formatting.
#   3. Remove the now-unwanted blank "No Spacing" paragraph that followed
#      "...Board of Directors".

$d = $word.ActiveDocument

# 1) Update the letter date (September 19 -> September 21, 2025), rewriting
#    only the matched run so the paragraph's bookmark/pPr stay untouched.
$searchRange = $d.Content
$dateFound = $searchRange.Find.Execute("September 19, 2025")
if ($dateFound) {
    $dateRange = $d.Range($searchRange.Start, $searchRange.End)
    $dateXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
      '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
      '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
      '<w:p><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">September 21, 2025</w:t></w:r></w:p>' + `
      '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $dateRange.InsertXML($dateXml) | Out-Null
}

# 2) Split the mailing address paragraph ("929 Story Road, San Jose CA
#    95122") into two separate paragraphs, leaving the identical address
#    that lives in the "PROPERTY ADDRESS" table untouched.
$addressRange = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -eq "929 Story Road, San Jose CA 95122`r" -and `
        $candidate.Range.Information(12) -eq $false) {
        # Information(12) = wdWithInTable; skip the table copy of this text.
        $addressRange = $candidate.Range
        break
    }
}

if ($addressRange -ne $null) {
    $addressXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
      '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
      '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
      '<w:p><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' + `
      '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">929 Story Road</w:t></w:r></w:p>' + `
      '<w:p><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' + `
      '<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">San Jose, CA 95122</w:t></w:r></w:p>' + `
      '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $addressRange.InsertXML($addressXml) | Out-Null
}

# 3) Drop the blank "No Spacing" paragraph right after "...Board of Directors".
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -match "Board of Directors") {
        $blank = $d.Paragraphs.Item($i + 1)
        if ($blank.Range.Text -eq "`r" -and $blank.Style.NameLocal -eq "No Spacing") {
            $blank.Range.Delete() | Out-Null
        }
        break
    }
}
